$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    We build it by copy/pasting the existing bold paragraph near the end of
#    the document (which already has the same leading-empty-run + bold-run
#    structure we want), so the new paragraph keeps that same run layout,
#    then we edit the bold run's text in place and append the rest as a
#    plain (non-bold) run.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$boldSourcePara = $d.Paragraphs(55)   # "Play Beer Bonanza Free Slot - Review & Demo Game" (bold)

$boldSourcePara.Range.Copy()

$afterTitle = $titlePara.Range
$afterTitle.Collapse(0)
$afterTitle.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

# Rewrite the bold run's text in place (keeps the leading empty run + bold
# run/rPr intact instead of collapsing the paragraph's runs).
$metaPara = $d.Paragraphs(2)
$metaFull = $metaPara.Range
$metaBoldRange = $d.Range($metaFull.Start, $metaFull.End - 1)
$metaBoldRange.Text = "Meta description"

# Append the (non-bold) rest of the sentence right after the bold run.
$metaPara = $d.Paragraphs(2)
$afterBoldPos = $metaPara.Range.End - 1
$afterBoldRange = $d.Range($afterBoldPos, $afterBoldPos)
$afterBoldRange.InsertAfter(": Read our unbiased review of Beer Bonanza, an Oktoberfest-themed slot game from BGaming. Play it for free or with real money at select casinos.")

# ---------------------------------------------------------------------------
# 2) Remove the old bold "Play Beer Bonanza..." paragraph near the end, and
#    replace the italic paragraph's text with the new image-prompt text.
# ---------------------------------------------------------------------------

$boldPara = $d.Paragraphs(56)   # shifted by +1 because of the paragraph inserted in step 1
$boldPara.Range.Delete()

$italicPara = $d.Paragraphs(56)
$italicFull = $italicPara.Range
$italicTextRange = $d.Range($italicFull.Start, $italicFull.End - 1)
$italicTextRange.Text = "Please create a feature image fitting the game ""Beer Bonanza"" with the following specifications: - Cartoon style - Happy Maya warrior with glasses"

Write-Output "done"
